# Update countries & provincias Spain
# - Reorder "Trinidad yTobago" ahead of "Ruanda" (rows 125/126) with refreshed case numbers
# - Reorder "Guyana" ahead of "Bahamas" / "San Martin (Parte Holandesa)" (rows 153-155) with refreshed case numbers
# - Refresh the "Datos actualizados..." timestamp

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 125-126: Trinidad yTobago now outranks Ruanda ---
$ws.Cells.Item(125, 1).Value2 = "Trinidad yTobago"
$ws.Cells.Item(125, 2).Value2 = 105
$ws.Cells.Item(125, 3).Value2 = 1
$ws.Cells.Item(125, 4).Value2 = 1
$ws.Cells.Item(125, 5).Value2 = 97
$ws.Cells.Item(125, 6).Value2 = 0
$ws.Cells.Item(125, 7).Value2 = 0
$ws.Cells.Item(125, 8).Value2 = 7

$ws.Cells.Item(126, 1).Value2 = "Ruanda"
$ws.Cells.Item(126, 2).Value2 = 104
$ws.Cells.Item(126, 3).Value2 = 0
$ws.Cells.Item(126, 4).Value2 = 4
$ws.Cells.Item(126, 5).Value2 = 100
$ws.Cells.Item(126, 6).Value2 = 0
$ws.Cells.Item(126, 7).Value2 = 0
$ws.Cells.Item(126, 8).Value2 = 0

# --- Rows 153-155: Guyana now outranks Bahamas / San Martin (Parte Holandesa) ---
$ws.Cells.Item(153, 1).Value2 = "Guyana"
$ws.Cells.Item(153, 2).Value2 = 29
$ws.Cells.Item(153, 3).Value2 = 5
$ws.Cells.Item(153, 4).Value2 = 0
$ws.Cells.Item(153, 5).Value2 = 25
$ws.Cells.Item(153, 6).Value2 = 0
$ws.Cells.Item(153, 7).Value2 = 0
$ws.Cells.Item(153, 8).Value2 = 4

$ws.Cells.Item(154, 1).Value2 = "Bahamas"
$ws.Cells.Item(154, 2).Value2 = 29
$ws.Cells.Item(154, 3).Value2 = 0
$ws.Cells.Item(154, 4).Value2 = 4
$ws.Cells.Item(154, 5).Value2 = 20
$ws.Cells.Item(154, 6).Value2 = 1
$ws.Cells.Item(154, 7).Value2 = 0
$ws.Cells.Item(154, 8).Value2 = 5

$ws.Cells.Item(155, 1).Value2 = "San Martin (Parte Holandesa)"
$ws.Cells.Item(155, 2).Value2 = 25
$ws.Cells.Item(155, 3).Value2 = 0
$ws.Cells.Item(155, 4).Value2 = 1
$ws.Cells.Item(155, 5).Value2 = 20
$ws.Cells.Item(155, 6).Value2 = 0
$ws.Cells.Item(155, 7).Value2 = 0
$ws.Cells.Item(155, 8).Value2 = 4

# --- Refresh "updated at" timestamp ---
$ws.Range("A1").Value2 = "Datos actualizados a 6 de Abril de 2020 a las 04:22"
